$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: remove the "Will Cite" (G5) mark for the Bracket Signatures entry
$ws.Range("G5").ClearContents()

# Row 12: add a (blank-ish/space) Year value for the Chen & Hwang entry
$ws.Range("D12").Value = " "

# New row 20: "Do Stronger Players Win More Knockout Tournaments" reference
# (order matters here so new shared-string indices line up with the source
#  workbook: Link, then Title, then Description, after D12's space string)
$ws.Range("E20").Value = "https://www-jstor-org.ezp-prod1.hul.harvard.edu/stable/pdf/2286606.pdf?refreqid=excelsior%3Ac6f9f27cd7701f5e7471bf45a00f7c81&ab_segments=&origin=&initiator=&acceptTC=1"
$ws.Range("B20").Value = "Do Stronger Players Win More Knockout Tournaments"
$ws.Range("H20").Value = "defined ordered"
$ws.Range("C20").Value = "Robert Chen and F. K. Hwang"
$ws.Range("D20").Value = 1978
$ws.Range("F20").Value = "Ordered Brackets"
$ws.Range("G20").Value = "x"

# Update the active selection to reflect where the author left off editing
$ws.Range("C6").Select() | Out-Null

# Set the sheet up for portrait-oriented printing
$ws.PageSetup.Orientation = 1
